$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.280.49'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '3.592.71'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '578.59'
$ws.Range('E5').Value = '  -2.57%  '
$ws.Range('D6').Value = '190.88'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('D8').Value = '3.585.57'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '0.178'
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('D12').Value = '56.66'
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').Value = '9.76'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = '4.171.55'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '20.18'
$ws.Range('E16').Value = '  +4.43%  '
$ws.Range('D17').Value = '3.591.91'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '70.175.04'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = '12.54'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '474.50'
$ws.Range('E22').Value = '  -5.24%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = '19.46'
$ws.Range('E23').Value = '  +12.89%  '
$ws.Range('D24').Value = '5.12'
$ws.Range('E24').Value = '  -6.69%  '
$ws.Range('D25').Value = '4.38'
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('D26').Value = '88.76'
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').Value = '11.12'
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('D29').Value = '9.28'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').Value = '7.75'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = '32.16'
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').Value = '0.120'
$ws.Range('E32').Value = '  +4.07%  '
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').Value = '66.12'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').Value = '588.48'
$ws.Range('E35').Value = '  -3.95%  '
$ws.Range('D36').Value = '39.64'
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  -3.72%  '
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('D41').Value = '3.56'
$ws.Range('E41').Value = '  -2.22%  '
$ws.Range('D42').Value = '2.93'
$ws.Range('E42').Value = '  +7.78%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.231.44'
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '3.15'
$ws.Range('E44').Value = '  +8.25%  '
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D47').Value = '9.66'
$ws.Range('E47').Value = '  +5.64%  '
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('E51').Value = '  -2.41%  '
